$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The ramp-chart's highlighted rectangle ("Rectangle 32") was resized/moved:
# it now starts a bit higher and is noticeably taller, to properly surround
# the highlighted ramp rows.
$sh = $s.Shapes.Item("Rectangle 32")
$sh.Left = 120.96004
$sh.Top = 163.7554
$sh.Width = 341.28004
$sh.Height = 34.8009

# The three arrow connectors that point from the "Ramps" label down into the
# highlighted area were stretched to match the new (taller) highlight box.
$sh = $s.Shapes.Item("Straight Arrow Connector 38")
$sh.Left = 350.88
$sh.Top = 158.4903
$sh.Width = 28.44
$sh.Height = 27.73405

$sh = $s.Shapes.Item("Straight Arrow Connector 40")
$sh.Left = 388.32
$sh.Top = 158.4903
$sh.Width = 0
$sh.Height = 27.734

$sh = $s.Shapes.Item("Straight Arrow Connector 42")
$sh.Left = 396.2401
$sh.Top = 158.4903
$sh.Width = 25.67996
$sh.Height = 27.1879

# The "Ramps" label itself shifted slightly left/up to stay centred above the
# (now wider-looking) arrow fan-out.
$sh = $s.Shapes.Item("TextBox 43")
$sh.Left = 365.91
$sh.Top = 140.9525
